$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A21").Value = "Nbr_Sign_Life"
$ws.Range("B21").Value = "entier"
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = "Nombre de zombies qui peuvent passer sur un panneau avant de la casser"

$ws.Range("A21:C21").VerticalAlignment = -4160

$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B21").Select()
